$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Semana 6 de 2026: refresh Poisson table -----------------------------
# Two new notifiable-event rows appear in this week's extract (codes 420
# and 740); every other row keeps its event code/name and only the
# Esperado/Observado/valor p numbers are refreshed.

# Insert rows at their FINAL target position, in ascending order, so each
# insert only has to push down the rows still below it.
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(27).Insert()

# Column A stores event codes as text (e.g. "113"); force text format on
# the new cells so COM does not silently coerce them to numbers.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A27").NumberFormat = "@"

$ws.Range("A18").Value = "420"
$ws.Range("B18").Value = "Leishmaniasis cutanea"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 1

$ws.Range("A27").Value = "740"
$ws.Range("B27").Value = "Sifilis congenita"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0

# Drop the explicit style the NumberFormat tweak above left behind, so the
# new cells match the unstyled (default) look of the rest of column A.
$ws.Range("A18").Style = "Normal"
$ws.Range("A27").Style = "Normal"

# --- Refresh Esperado / Observado / valor p for the remaining events -----
$updates = @(
  ,@(2, 2, 3, 0.18)
  ,@(3, 0, 0, 1)
  ,@(4, 5, 8, 0.07000000000000001)
  ,@(5, 2, 5, 0.04)
  ,@(6, 2, 2, 0.27)
  ,@(7, 46, 71, 0)
  ,@(8, 0, 0, 1)
  ,@(9, 1, 0, 0.37)
  ,@(10, 2, 5, 0.04)
  ,@(11, 67, 0, 0)
  ,@(12, 0, 1, 0)
  ,@(13, 1, 0, 0.37)
  ,@(15, 9, 8, 0.13)
  ,@(16, 1, 0, 0.37)
  ,@(17, 8, 2, 0.01)
  ,@(19, 1, 2, 0.18)
  ,@(22, 7, 7, 0.15)
  ,@(26, 1, 0, 0.37)
  ,@(28, 2, 1, 0.27)
  ,@(30, 7, 6, 0.15)
  ,@(31, 8, 5, 0.09)
  ,@(32, 8, 7, 0.14)
)
foreach ($u in $updates) {
  $ws.Cells.Item($u[0], 3).Value = $u[1]
  $ws.Cells.Item($u[0], 4).Value = $u[2]
  $ws.Cells.Item($u[0], 5).Value = $u[3]
}
